$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the touched Price/Volume cells to remain plain text (matching the
# original inlineStr cells) instead of being auto-parsed into numbers/percentages.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "332.15"
$ws.Range("E2").Value = "1.62%"
$ws.Range("E3").Value = "3.56%"
$ws.Range("D4").Value = "5.699"
$ws.Range("E4").Value = "3.43%"
$ws.Range("D5").Value = "0.08371"
$ws.Range("E5").Value = "4.38%"
$ws.Range("D6").Value = "2.038"
$ws.Range("E6").Value = "2.61%"
$ws.Range("D7").Value = "0.9812"
$ws.Range("E7").Value = "3.47%"
$ws.Range("E8").Value = "0.94%"
$ws.Range("D9").Value = "0.1166"
$ws.Range("E9").Value = "1.64%"
$ws.Range("D10").Value = "0.1942"
$ws.Range("E10").Value = "5.75%"
$ws.Range("D11").Value = "10.36"
$ws.Range("E11").Value = "-14.71%"
$ws.Range("D12").Value = "0.1009"
$ws.Range("E12").Value = "2.95%"
$ws.Range("D13").Value = "0.04668"
$ws.Range("E13").Value = "1.23%"
$ws.Range("D14").Value = "0.1058"
$ws.Range("E14").Value = "-0.74%"
$ws.Range("D15").Value = "0.001290"
$ws.Range("E15").Value = "1.78%"
$ws.Range("D16").Value = "0.006057"
$ws.Range("E16").Value = "5.87%"
$ws.Range("D17").Value = "3.370"
$ws.Range("E17").Value = "0.11%"
$ws.Range("D18").Value = "4.465"
$ws.Range("E18").Value = "4.14%"
$ws.Range("E20").Value = "-0.39%"
$ws.Range("D21").Value = "0.2593"
$ws.Range("E21").Value = "1.96%"
$ws.Range("D22").Value = "0.04210"
$ws.Range("E22").Value = "3.17%"
$ws.Range("E23").Value = "5.27%"
$ws.Range("D24").Value = "0.004589"
$ws.Range("E24").Value = "6.23%"
$ws.Range("E25").Value = "7.67%"
$ws.Range("E26").Value = "0.00%"
$ws.Range("D38").Value = "0.02778"
$ws.Range("E38").Value = "8.34%"
$ws.Range("D39").Value = "0.05818"
$ws.Range("E39").Value = "4.88%"
$ws.Range("D40").Value = "0.007739"
$ws.Range("E40").Value = "2.73%"
$ws.Range("E41").Value = "3.16%"
$ws.Range("D42").Value = "0.007195"
$ws.Range("E42").Value = "-5.38%"
$ws.Range("E43").Value = "-1.98%"
$ws.Range("D44").Value = "0.008180"
$ws.Range("E44").Value = "-3.94%"
$ws.Range("E45").Value = "1.19%"
$ws.Range("E46").Value = "0.10%"
$ws.Range("E47").Value = "-0.15%"
$ws.Range("D48").Value = "0.003492"
$ws.Range("E48").Value = "89.27%"
$ws.Range("E49").Value = "-0.80%"
$ws.Range("E50").Value = "0.10%"
$ws.Range("E51").Value = "0.10%"

# Restore the original (default) cell style now that the text values are locked in.
$ws.Range("D2:E51").Style = "Normal"
